$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tall "Маловато будет =(\n#ДеньРадио" entry currently lives in A8 (the row
# with the large 71.25 row height). The data loader fix moves this entry to the
# end of the list (A15) and leaves A8 blank (but keeps its formatting/height).

$sourceCell = $ws.Range("A8")
$movedText = $sourceCell.Value2

$targetCell = $ws.Range("A15")

# Reset A15 to the plain/default style first so that applying WrapText below
# reuses the existing "wrap text" cell style (same one used by A8) instead of
# merging with whatever formatting A15 already had.
$targetCell.Style = "Normal"
$targetCell.Value = $movedText
$targetCell.WrapText = $true
$ws.Rows.Item(15).RowHeight = 71.25

# Clear the old location; row height/style for A8 stay as-is.
$sourceCell.ClearContents()
